$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.507.14"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.958.58"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'244.17"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "'58.65"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("E10").Value = "  -5.57%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "'22.09"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "'0.829"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "2.245.32"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'13.71"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "'5.28"
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "1.958.69"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "36.459.02"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'69.71"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "'228.36"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +2.10%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "'0.138"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").Value = "'160.32"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "'19.44"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "'0.120"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("D32").Value = "'4.71"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "'0.0620"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("D34").Value = "'4.31"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "'2.25"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "  +10.20%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").Value = "'5.79"
$ws.Range("E39").Value = "  -9.58%  "
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'16.06"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "1.366.38"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").Value = "'87.88"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").Value = "'7.14"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("D49").Value = "'2.83"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "2.135.78"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E51").Value = "  -5.13%  "
